# Update Beta group significance (Family/Tribe) table with refreshed
# permutation-test statistics (pseudo-F, p-value, q-value) and corrected
# group-pair ordering/sample sizes for several rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 5.547265439098457
$ws.Range("E3").Value = 5.235264861421014
$ws.Range("E4").Value = 3.933351878520496
$ws.Range("E5").Value = 2.163743604272656
$ws.Range("F5").Value = 0.01
$ws.Range("G5").Value = 0.025
$ws.Range("E6").Value = 1.6952695747238
$ws.Range("F6").Value = 0.053
$ws.Range("G6").Value = 0.08833333333333332
$ws.Range("E7").Value = 2.556227569289787
$ws.Range("F7").Value = 0.049
$ws.Range("G7").Value = 0.08833333333333332
$ws.Range("E8").Value = 1.390499738681556
$ws.Range("F8").Value = 0.135
$ws.Range("G8").Value = 0.1928571428571429
$ws.Range("A9").Value = 'Mugilidae'
$ws.Range("C9").Value = 13
$ws.Range("E9").Value = 1.018975555051429
$ws.Range("F9").Value = 0.456
$ws.Range("G9").Value = 0.57
$ws.Range("A10").Value = 'Haplochrominae'
$ws.Range("C10").Value = 15
$ws.Range("E10").Value = 0.9148717685693164
$ws.Range("F10").Value = 0.555
$ws.Range("G10").Value = 0.6166666666666667
$ws.Range("E12").Value = 4.123448265991692
$ws.Range("G12").Value = 0.0025
$ws.Range("B13").Value = 'Mugilidae'
$ws.Range("C13").Value = 67
$ws.Range("E13").Value = 2.717208054916494
$ws.Range("G13").Value = 0.0025
$ws.Range("B14").Value = 'Nemacheilidae'
$ws.Range("C14").Value = 85
$ws.Range("E14").Value = 3.295176388509694
$ws.Range("G14").Value = 0.0025
$ws.Range("A15").Value = 'Cyprinidae'
$ws.Range("B15").Value = 'Tilapiinae'
$ws.Range("C15").Value = 86
$ws.Range("E15").Value = 5.620056146322516
$ws.Range("G15").Value = 0.0025
$ws.Range("A16").Value = 'Haplochrominae'
$ws.Range("B16").Value = 'Nemacheilidae'
$ws.Range("C16").Value = 29
$ws.Range("E16").Value = 4.063042248649333
$ws.Range("G16").Value = 0.0025
$ws.Range("A17").Value = 'Nemacheilidae'
$ws.Range("B17").Value = 'Tilapiinae'
$ws.Range("C17").Value = 43
$ws.Range("E17").Value = 4.55369472125412
$ws.Range("F17").Value = 0.001
$ws.Range("G17").Value = 0.0025
$ws.Range("E18").Value = 2.178984864939341
$ws.Range("E19").Value = 2.871351695669131
$ws.Range("F19").Value = 0.003
$ws.Range("G19").Value = 0.005625
$ws.Range("E20").Value = 2.312906014824591
$ws.Range("F20").Value = 0.021
$ws.Range("G20").Value = 0.035
$ws.Range("E21").Value = 1.433109065074006
$ws.Range("F21").Value = 0.095
$ws.Range("G21").Value = 0.1425
$ws.Range("A22").Value = 'Mugilidae'
$ws.Range("B22").Value = 'Tilapiinae'
$ws.Range("C22").Value = 25
$ws.Range("E22").Value = 1.139804165827504
$ws.Range("F22").Value = 0.252
$ws.Range("G22").Value = 0.3436363636363637
$ws.Range("A23").Value = 'Haplochrominae'
$ws.Range("B23").Value = 'Poeciliidae'
$ws.Range("C23").Value = 10
$ws.Range("E23").Value = 1.05697324977583
$ws.Range("F23").Value = 0.317
$ws.Range("G23").Value = 0.3657692307692307
$ws.Range("A24").Value = 'Poeciliidae'
$ws.Range("C24").Value = 24
$ws.Range("E24").Value = 1.083863234109317
$ws.Range("F24").Value = 0.315
$ws.Range("G24").Value = 0.3657692307692307
$ws.Range("E25").Value = 0.9921702742492496
$ws.Range("F25").Value = 0.406
$ws.Range("G25").Value = 0.435
$ws.Range("E26").Value = 0.8031809414898595
$ws.Range("F26").Value = 0.6970000000000001
$ws.Range("G26").Value = 0.6970000000000001
